# Cap Ret per Unit Net Loss.xlsx
# Enable economic retirements of coal plants due to ETS by requiring plants
# to return a profit in SoFCtMbCtPR: bump the "natural gas combined cycle"
# capacity-retired-per-unit-net-loss parameter on the CRpUNL sheet from
# 0.03 to 0.035 MW/($/MW).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRpUNL")

# Row 4 = "natural gas combined cycle" -> column B holds the MW-retired slope.
$ws.Range("B4").Value = 0.035

# Bring the CRpUNL sheet to the front / make it the active tab+selection,
# matching the saved view state captured in the workbook.
$ws.Activate()
$ws.Range("B5").Select()
